$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.577122092247009
$ws.Range("B1").Value = 2.665462493896484
$ws.Range("C1").Value = 3.001894950866699
$ws.Range("D1").Value = 2.820022344589233
$ws.Range("E1").Value = 3.233738422393799
